# Update (Analyze PO & Forecast)
# Revises the MyForecast figures on "Forecast Comparison" and the
# corresponding rolled-up totals on "Summary".

$wb = $excel.ActiveWorkbook

# --- Sheet: Forecast Comparison -------------------------------------------
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")

$wsForecast.Range("D2").Value = 102
$wsForecast.Range("D3").Value = 104
$wsForecast.Range("D4").Value = 113
$wsForecast.Range("D12").Value = 124
$wsForecast.Range("D13").Value = 117
$wsForecast.Range("D14").Value = 119
$wsForecast.Range("D17").Value = 120

# --- Sheet: Summary ---------------------------------------------------------
# These metrics are stored as text (not numbers), so force a text number
# format before writing the digits, then clear the format again so the
# cell's style index is left untouched (matches the source formatting).
$wsSummary = $wb.Worksheets.Item("Summary")

foreach ($pair in @(
        @{ Cell = "B9";  Text = "1982" },
        @{ Cell = "B10"; Text = "986" },
        @{ Cell = "B11"; Text = "449" },
        @{ Cell = "B14"; Text = "102" }
    )) {
    $cell = $wsSummary.Range($pair.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $pair.Text
    $cell.ClearFormats()
}
